$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.452.08'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '3.500.01'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.58'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.95'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('D7').Value = '3.496.87'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.487'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.56'
$ws.Range('E11').Value = '  +7.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.45'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '4.093.19'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '3.493.49'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').Value = '67.370.15'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('E19').Value = '  +2.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.59'
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('E21').Value = '  +5.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '445.50'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.632'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').Value = '3.641.39'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('E28').Value = '  +4.92%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.07'
$ws.Range('E29').Value = '  -2.29%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.65'
$ws.Range('E31').Value = '  +6.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.169'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.01'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.66'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').Value = '3.495.21'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.32'
$ws.Range('E40').Value = '  +6.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '175.01'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0891'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.47'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.886'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.99'
$ws.Range('E46').Value = '  +7.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.13'
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.30'
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('E49').Value = '  -3.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.995'
$ws.Range('E51').Value = '  -0.75%  '
